{"js": "// The commit removes the trailing \"signature block\" paragraphs at the very\n// end of the document body: four empty spacer paragraphs, a paragraph with\n// nine tabs + the date (\"1 Desember 2021, 07:03\"), another empty paragraph,\n// and the final tabs/spaces + \"Baharudin Pratama\" paragraph. After the\n// change, the document body's last paragraph is the one holding the final\n// inline picture (Picture 13), immediately followed by the section\n// properties (<w:sectPr>).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst total = paragraphs.items.length;\nconst numToRemove = 6; // the 6 trailing paragraphs after the last picture\n\n// Delete from a single, already-loaded collection so the shim's live-proxy /\n// block-shift bookkeeping (mirroring real Word/Office.js semantics) keeps\n// every item's identity correct as earlier deletes shift later indices.\nfor (let i = total - numToRemove; i < total; i++) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# The commit removes the trailing \"signature block\" paragraphs at the very\n# end of the document body: four empty spacer paragraphs, a paragraph with\n# nine tabs + the date (\"1 Desember 2021, 07:03\"), another empty paragraph,\n# and the final tabs/spaces + \"Baharudin Pratama\" paragraph. After the\n# change, the document body's last paragraph is the one holding the final\n# inline picture (Picture 13), immediately followed by the section\n# properties.\n$d = $word.ActiveDocument\n\n$numToRemove = 6\n$count = $d.Paragraphs.Count\n$firstToRemove = $count - $numToRemove + 1   # 1-indexed\n\n# Build one Range spanning from the start of the first paragraph to remove\n# through the very end of the document's main story, then delete it in a\n# single operation (avoids churn from re-resolving \"the last paragraph\"\n# after each individual delete).\n$startPos = $d.Paragraphs.Item($firstToRemove).Range.Start\n$endPos = $d.Content.End\n$r = $d.Range($startPos, $endPos)\n$r.Delete()\n"}
